# Commit: "Codebase should now work with command prompt starts"
#
# The workbook originally had two sheets:
#   1) "Sheet 1 - 20230603-a1r-yc-sessi"  (the transcript: speaker/time/text)
#   2) "proposal"                          (a simple list of proposals)
#
# The edit removes the "proposal" sheet entirely and renames the remaining
# transcript sheet to "in" (its data/columns/values are unchanged).

$wb = $excel.ActiveWorkbook
[void]($excel.DisplayAlerts = $false)

# Remove the "proposal" worksheet if present.
foreach ($sheet in @($wb.Worksheets)) {
    if ($sheet.Name -eq "proposal") {
        [void]$sheet.Delete()
    }
}

# Rename the remaining (first/only) worksheet to "in".
$ws = $wb.Worksheets.Item(1)
[void]($ws.Name = "in")
